# New crime data collected
# Updates the CompStat weekly report: header (commissioner name, volume/issue
# number, reporting week dates) and the crime-complaint statistics table
# (rows 15-30) with the newly collected figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header text updates (shared-string backed cells)
# ---------------------------------------------------------------------------
$ws.Range("M6").Value = "Edward A. Caban"
$ws.Range("A8").Value = "Volume 30   Number  27"
$ws.Range("C8").Value = "Report Covering the Week  7/3/2023  Through  7/9/2023"

# ---------------------------------------------------------------------------
# Helper: formats used by the numeric columns in the crime-complaints table.
#   "#,##0"                      -> plain integer / count column style
#   "#,##0.0;""-""#,##0.0"       -> percent-change column style
# ---------------------------------------------------------------------------
$intFmt = "#,##0"
$pctFmt = "#,##0.0;""-""#,##0.0"

# A cell that already carries the plain "General" text style (s=14) used by
# the "0" / "***.*" placeholder cells throughout the table; copying its
# format onto a cell (after forcing the cell to Text so the literal isn't
# re-interpreted as a number) reproduces that exact style.
$blankStyleSource = $ws.Range("A14")

function Set-PlaceholderText($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $blankStyleSource.Copy()
    $rng.PasteSpecial(-4122)
}

function Set-IntValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = $intFmt
    $rng.Value = $value
}

function Set-PctValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = $pctFmt
    $rng.Value = $value
}

# ---------------------------------------------------------------------------
# Row 15 (Rape)
# ---------------------------------------------------------------------------
$ws.Range("L15").Value = 14.285714285714
$ws.Range("N15").Value = -33.333333333333

# ---------------------------------------------------------------------------
# Row 16 (Robbery)
# ---------------------------------------------------------------------------
$ws.Range("D16").Value = 5
$ws.Range("E16").Value = -60
$ws.Range("F16").Value = 8
$ws.Range("G16").Value = 13
$ws.Range("H16").Value = -38.461538461538
$ws.Range("I16").Value = 72
$ws.Range("J16").Value = 61
$ws.Range("K16").Value = 18.032786885245
$ws.Range("L16").Value = 35.849056603773
$ws.Range("M16").Value = -13.253012048192
$ws.Range("N16").Value = -82.524271844660

# ---------------------------------------------------------------------------
# Row 17 (Fel. Assault)
# ---------------------------------------------------------------------------
$ws.Range("C17").Value = 4
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = -20
$ws.Range("F17").Value = 15
$ws.Range("G17").Value = 17
$ws.Range("H17").Value = -11.764705882352
$ws.Range("I17").Value = 111
$ws.Range("J17").Value = 86
$ws.Range("K17").Value = 29.069767441860
$ws.Range("L17").Value = 18.085106382978
$ws.Range("M17").Value = 98.214285714285
$ws.Range("N17").Value = -29.746835443038

# ---------------------------------------------------------------------------
# Row 18 (Burglary) -- D18/E18 switch from "N/A" placeholders to real numbers
# ---------------------------------------------------------------------------
$ws.Range("C18").Value = 3
Set-IntValue "D18" 1
Set-PctValue "E18" 200
$ws.Range("F18").Value = 9
$ws.Range("G18").Value = 6
$ws.Range("H18").Value = 50
$ws.Range("I18").Value = 116
$ws.Range("J18").Value = 88
$ws.Range("K18").Value = 31.818181818181
$ws.Range("L18").Value = 52.631578947368
$ws.Range("M18").Value = -20.547945205479
$ws.Range("N18").Value = -84.533333333333

# ---------------------------------------------------------------------------
# Row 19 (Gr. Larceny)
# ---------------------------------------------------------------------------
$ws.Range("C19").Value = 16
$ws.Range("D19").Value = 13
$ws.Range("E19").Value = 23.076923076923
$ws.Range("F19").Value = 54
$ws.Range("G19").Value = 55
$ws.Range("H19").Value = -1.818181818181
$ws.Range("I19").Value = 337
$ws.Range("J19").Value = 365
$ws.Range("K19").Value = -7.671232876712
$ws.Range("L19").Value = 17.832167832167
$ws.Range("M19").Value = 52.488687782805
$ws.Range("N19").Value = -13.810741687979

# ---------------------------------------------------------------------------
# Row 20 (G.L.A.)
# ---------------------------------------------------------------------------
$ws.Range("C20").Value = 3
$ws.Range("D20").Value = 5
$ws.Range("E20").Value = -40
$ws.Range("G20").Value = 21
$ws.Range("H20").Value = -47.619047619047
$ws.Range("I20").Value = 72
$ws.Range("J20").Value = 71
$ws.Range("K20").Value = 1.408450704225
$ws.Range("L20").Value = 84.615384615384
$ws.Range("M20").Value = -18.181818181818
$ws.Range("N20").Value = -92.258064516129

# ---------------------------------------------------------------------------
# Row 21 (TOTAL)
# ---------------------------------------------------------------------------
$ws.Range("C21").Value = 28
$ws.Range("E21").Value = -3.448275862068
$ws.Range("F21").Value = 97
$ws.Range("G21").Value = 114
$ws.Range("H21").Value = -14.912280701754
$ws.Range("I21").Value = 717
$ws.Range("J21").Value = 685
$ws.Range("K21").Value = 4.671532846715
$ws.Range("L21").Value = 28.956834532374
$ws.Range("M21").Value = 18.708609271523
$ws.Range("N21").Value = -73.004518072289

# ---------------------------------------------------------------------------
# Row 22 (Transit) -- D22/E22 switch from real numbers to "N/A" placeholders
# ---------------------------------------------------------------------------
Set-PlaceholderText "D22" "0"
Set-PlaceholderText "E22" "***.*"
$ws.Range("F22").Value = 3
$ws.Range("H22").Value = 200
$ws.Range("I22").Value = 7
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 250
$ws.Range("M22").Value = -61.111111111111

# ---------------------------------------------------------------------------
# Row 24 (Petit Larceny)
# ---------------------------------------------------------------------------
$ws.Range("C24").Value = 24
$ws.Range("D24").Value = 45
$ws.Range("E24").Value = -46.666666666666
$ws.Range("F24").Value = 151
$ws.Range("G24").Value = 176
$ws.Range("H24").Value = -14.204545454545
$ws.Range("I24").Value = 948
$ws.Range("J24").Value = 886
$ws.Range("K24").Value = 6.997742663656
$ws.Range("L24").Value = 57.213930348258
$ws.Range("M24").Value = 97.5

# ---------------------------------------------------------------------------
# Row 25 (Misd. Assault)
# ---------------------------------------------------------------------------
$ws.Range("C25").Value = 15
$ws.Range("D25").Value = 10
$ws.Range("E25").Value = 50
$ws.Range("G25").Value = 38
$ws.Range("H25").Value = 39.473684210526
$ws.Range("I25").Value = 318
$ws.Range("J25").Value = 228
$ws.Range("K25").Value = 39.473684210526
$ws.Range("L25").Value = 45.205479452054
$ws.Range("M25").Value = 27.2

# ---------------------------------------------------------------------------
# Row 26 (UCR Rape*)
# ---------------------------------------------------------------------------
$ws.Range("L26").Value = 57.142857142857

# ---------------------------------------------------------------------------
# Row 27 (Other Sex Crimes) -- C27/G27/H27 switch to "N/A" placeholders
# ---------------------------------------------------------------------------
Set-PlaceholderText "C27" "0"
Set-PlaceholderText "G27" "0"
Set-PlaceholderText "H27" "***.*"

# ---------------------------------------------------------------------------
# Row 28 (Shooting Vic.) -- D28/E28/G28/H28 switch from placeholders to numbers
# ---------------------------------------------------------------------------
Set-IntValue "D28" 1
Set-PctValue "E28" -100
Set-IntValue "G28" 1
Set-PctValue "H28" -100
$ws.Range("J28").Value = 2
$ws.Range("K28").Value = -50

# ---------------------------------------------------------------------------
# Row 29 (Shooting Inc.) -- D29/E29/G29/H29 switch from placeholders to numbers
# ---------------------------------------------------------------------------
Set-IntValue "D29" 1
Set-PctValue "E29" -100
Set-IntValue "G29" 1
Set-PctValue "H29" -100
$ws.Range("J29").Value = 2
$ws.Range("K29").Value = -50

# ---------------------------------------------------------------------------
# Row 30 (Hate Crimes) -- D30/E30 switch to "N/A" placeholders
# ---------------------------------------------------------------------------
Set-PlaceholderText "D30" "0"
Set-PlaceholderText "E30" "***.*"
